$d = $word.ActiveDocument

# Locate the target paragraph: "Answering the research questions for this topic ... grouping."
# It is the last paragraph in the document body (immediately before the sectPr).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Answering the research questions for this topic*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$r = $target.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">Answering the research questions for this topic requires demographic information about victims of the police, traffic stops metrics and county census data.  These sources can feed into a model that assesses the risk of mortality present to each class of citizens.  That assessment needs to express descriptive statistics across both horizontals (e.g., age and income) and verticals (e.g., nationality and race) groupings.  These groups can bubble up latent feature dependencies, such as poverty-stricken people might be </w:t></w:r><w:r><w:t xml:space="preserve">more influential than </w:t></w:r><w:r><w:t xml:space="preserve">race.  </w:t></w:r><w:r><w:t xml:space="preserve">Uncovering </w:t></w:r><w:r><w:t xml:space="preserve">these intricate details requires </w:t></w:r><w:r><w:t xml:space="preserve">further investigation into the </w:t></w:r><w:r><w:t xml:space="preserve">relative log-likelihood </w:t></w:r><w:r><w:t xml:space="preserve">between </w:t></w:r><w:r><w:t>grouping.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Another critical piece of the puzzle is determining the validity in defunding the police and pivoting toward civil service investments.  According to a cursory investigation, most efforts on this front have been symbolic at best</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="450063084"/><w:citation/></w:sdtPr><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Smi20 \l 1033 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Smith, 2020)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve">.  Instead, a model needs to exist for better </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>categorizing risk levels during police stops.  For instance, how many fatalities involve mental illness, drug abuse, or are unarmed?  These features might unlock additional latent features that enable the officer to operate differently while still ensuring personal safety.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Hypothesis Testing</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">A hypothesis </w:t></w:r></w:p>
'@

$r.InsertXML($xml)
